$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$oldText = $cellA1.Value2
$newText = $oldText -replace [regex]::Escape("1000 Bs = 3.4 = 13110.39 pesos"), "1000 Bs = 3.41 = 13116.26 pesos"
$newText = $newText -replace [regex]::Escape("13110.39 pesos = 3.38 = 946.43 Bs"), "13116.26 pesos = 3.39 = 948.49 Bs"
$cellA1.Value = $newText

# --- Sheet "tasas": update the numeric rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 293.3
$wsTasas.Range("O10").Value = 3847
$wsTasas.Range("N12").Value = 3872
